# "Added Week-4 & updated Week-3 PPT"
# The only content change to this deck (Week 3 - PPT.pptx) is the title
# label on slide 1 flipping from "WEEK-2" to "WEEK-3" (the rest of the
# commit's XML diff is just incidental namespace/attribute reordering
# produced by PowerPoint's own re-serialization, not a semantic edit).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "WEEK-2") {
                $target = $shp
                break
            }
        }
    }
}

if ($target -eq $null) {
    # Fallback: known position of the "WEEK-2" label shape.
    $target = $s.Shapes.Item(7)
}

$target.TextFrame.TextRange.Text = "WEEK-3"
